# Actualización desde MV -datos-
# Updates the "Tasas de colocación 2021 - Diaria" sheet with new daily
# rows for 16, 20, 21, 22, 23, 24, 27, 28 and 29-09-2021, and corrects
# the 15-09-2021 "30 a 89 días, US$" figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the existing 15-09-2021 row (row 179), column G ---
$ws.Cells.Item(179, 7).Value = 2.11

# --- Append the new daily rows (180-188) ---
# Columns: A Serie, B interbancaria, C 1-30d, D 30-89d, E 90d-1a,
#          F 30-89d UF, G 90d-1a UF, H 30-89d US$, I 90d-1a US$
$newRows = @(
    @{ Row = 180; Serie = "16-09-2021"; B = 0.13; C = 0.5600000000000001; D = 0.3;  E = 1.35; F = 3.69; G = 2.12; H = 2.36;               I = 2.54 },
    @{ Row = 181; Serie = "20-09-2021"; B = 0.13; C = 0.85;               D = 0.35; E = 0.92; F = $null; G = 2.12; H = 2.5;                I = 2.17 },
    @{ Row = 182; Serie = "21-09-2021"; B = 0.13; C = 0.67;               D = 0.32; E = 0.86; F = 4.26; G = 2.04; H = 0.8100000000000001; I = 1.46 },
    @{ Row = 183; Serie = "22-09-2021"; B = 0.13; C = 0.76;               D = 0.32; E = 0.62; F = 1.57; G = 2.07; H = 1.97;               I = 2.06 },
    @{ Row = 184; Serie = "23-09-2021"; B = 0.13; C = 0.55;               D = 0.28; E = 0.58; F = 1.79; G = 1.84; H = 2.31;               I = 2.5  },
    @{ Row = 185; Serie = "24-09-2021"; B = 0.13; C = 0.73;               D = 0.36; E = 0.74; F = 1.8;  G = 2.61; H = 1.4;                I = 2.28 },
    @{ Row = 186; Serie = "27-09-2021"; B = 0.13; C = 0.5600000000000001; D = 0.31; E = 1.1;  F = 2.95; G = 1.88; H = 2.4;                I = 2.73 },
    @{ Row = 187; Serie = "28-09-2021"; B = 0.13; C = 0.82;               D = 0.39; E = 1.03; F = 1.96; G = 2.1;  H = 1.34;               I = 2.29 },
    @{ Row = 188; Serie = "29-09-2021"; B = 0.13; C = 0.47;               D = 0.41; E = 0.57; F = 2.2;  G = 1.99; H = 0.74;               I = 2.32 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Serie
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    if ($r.F -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r.F
    }
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}
